# Daily attendance processing - reorder the "Recorded By" (column G) list
# so the automation/system account name no longer leads the human
# recorder's name/email. For each data row, the comma-separated list of
# recorders is reversed, unless it already ends in "System" (which means
# the list is already in the desired, human-first order).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Text

    if ([string]::IsNullOrEmpty($v)) { continue }

    $parts = $v -split ", "
    $n = $parts.Length

    if ($n -lt 2) { continue }
    if ($parts[$n - 1] -eq "System") { continue }

    $rev = $parts[($n - 1)..0]
    $joined = [string]::Join(", ", $rev)

    $cell.Value = $joined
}
